# Table 19: "Added Above TAC to quotas Factor"
# - Header: SPECIES -> COUNTRIES
# - Header: COUNTRIES -> QUOTAS (+ line break + "cod", italic)
# - Header: QUOTAS -> QUOTAS (unchanged) (+ line break + "hake", italic)
# - Data rows: col1 "Atlantic cod" -> country code (was col2)
#              col2 country code -> new "cod" quota value
#              col3 old quota value -> new "hake" quota value

$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

function Set-CellText($table, $row, $col, $newText) {
    $c = $table.Cell($row, $col).Range
    $full = $c.Text
    $trimLen = $full.Length - 2
    $sub = $d.Range($c.Start, $c.Start + $trimLen)
    $sub.Text = $newText
}

function Append-CellBreakText($table, $row, $col, $newText) {
    # Appends a manual line break (vertical tab char) followed by newText
    # to the end of the existing cell text. Must use the cell's own Range
    # object (which spans through the cell-end markers) collapsed to its
    # end -- using a freshly constructed sub-range here does not reliably
    # insert content in this runtime.
    $c = $table.Cell($row, $col).Range
    $c.Collapse(0)
    $c.InsertAfter([char]11 + $newText)
}

# ---- Header row ----
Set-CellText $t 1 1 "COUNTRIES"
Set-CellText $t 1 2 "QUOTAS"
Append-CellBreakText $t 1 2 "cod"
Append-CellBreakText $t 1 3 "hake"

# ---- Data rows ----
$rows = @(
    @{ r = 2;  country = "BE"; cod = "0.488"; hake = "0.494" },
    @{ r = 3;  country = "DK"; cod = "0.175"; hake = "0.529" },
    @{ r = 4;  country = "DE"; cod = "0.521"; hake = "0.490" },
    @{ r = 5;  country = "EE"; cod = "0.480"; hake = "0.489" },
    @{ r = 6;  country = "IE"; cod = "0.460"; hake = "0.515" },
    @{ r = 7;  country = "ES"; cod = "0.639"; hake = "0.780" },
    @{ r = 8;  country = "FR"; cod = "0.476"; hake = "0.965" },
    @{ r = 9;  country = "LV"; cod = "0.458"; hake = "0.489" },
    @{ r = 10; country = "LT"; cod = "0.467"; hake = "0.489" },
    @{ r = 11; country = "NL"; cod = "0.485"; hake = "0.493" },
    @{ r = 12; country = "PL"; cod = "0.502"; hake = "0.489" },
    @{ r = 13; country = "PT"; cod = "0.520"; hake = "0.533" },
    @{ r = 14; country = "FI"; cod = "0.460"; hake = "0.489" },
    @{ r = 15; country = "SE"; cod = "0.365"; hake = "0.490" }
)

foreach ($row in $rows) {
    Set-CellText $t $row.r 1 $row.country
    Set-CellText $t $row.r 2 $row.cod
    Set-CellText $t $row.r 3 $row.hake
}

Write-Host "Edit complete"
